$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.096.89"
$ws.Range("E2").Value = "  +0.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.751.01"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "235.35"
$ws.Range("E5").Value = "  +4.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5291"
$ws.Range("E7").Value = "  +2.61%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2800"
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06200"
$ws.Range("E9").Value = "  +1.75%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.745.15"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07176"
$ws.Range("E11").Value = "  +2.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.47"
$ws.Range("E12").Value = "  +1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6476"
$ws.Range("E13").Value = "  +2.07%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.639"
$ws.Range("E14").Value = "  +3.23%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.73"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9996"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9993"
$ws.Range("E17").Value = "  -0.12%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "25.999.01"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.72"
$ws.Range("E19").Value = "  +2.56%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000006748"
$ws.Range("E20").Value = "  +2.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.967.86"
$ws.Range("E21").Value = "  +0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.327"
$ws.Range("E22").Value = "  +6.41%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.734"
$ws.Range("E23").Value = "  +3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.251"
$ws.Range("E24").Value = "  +3.07%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "138.50"
$ws.Range("E25").Value = "  +0.61%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.509"
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("E27").Value = "  +2.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.811"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "104.82"
$ws.Range("E29").Value = "  +2.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08293"
$ws.Range("E30").Value = "  +0.28%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.809"
$ws.Range("E31").Value = "  +5.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.656"
$ws.Range("E32").Value = "  +7.29%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04581"
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("E34").Value = "  +0.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("E35").Value = "  +4.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6361"
$ws.Range("E36").Value = "  +6.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.711"
$ws.Range("E37").Value = "  +1.57%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01604"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.955"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "100.78"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.3944"
$ws.Range("E42").Value = "  +3.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7456"
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.036"
$ws.Range("E44").Value = "  +3.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1150"
$ws.Range("E45").Value = "  +4.50%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.378"
$ws.Range("E46").Value = "  +2.09%  "
$ws.Range("E47").Value = "  -2.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "31.24"
$ws.Range("E48").Value = "  +5.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "54.34"
$ws.Range("E49").Value = "  +4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.599"
$ws.Range("E50").Value = "  +1.95%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3461"
$ws.Range("E51").Value = "  +2.42%  "
